$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing all existing data down by one row
$ws.Rows.Item(1).Insert()

# Fill in the new header row (row 1) with the inner-header labels
$headers = @("Название", "Норма азота", "Норма фосфора", "Норма калия", "Культура", "Район", "Цена", "Описание", "Назначение")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Style the header row: bold font + grey fill, applied once on A1 and then
# copied (format only) across the rest of the row so only a single new
# style definition is produced.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Interior.Color = 0xEAEAEA
$a1.Copy()
$ws.Range("B1:I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(1).RowHeight = 37.5

# Rename the product rows (formerly "10/11/12 удобрение") to "16/17/18 удобрение"
$ws.Range("A2").Value = "16 удобрение"
$ws.Range("A3").Value = "17 удобрение"
$ws.Range("A4").Value = "18 удобрение"

# Append a new empty row (23) below the last existing blank row, keeping
# the same formatting as the rest of the placeholder rows (B:D, style 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp = -4162
$newRow = $lastRow + 1
$ws.Range("B" + $lastRow + ":D" + $lastRow).Copy()
$ws.Range("B" + $newRow + ":D" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to match target
$ws.Range("A11").Select()
